$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.859.42"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.73%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.805.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.41%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9994"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4995"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.78%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3884"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09435"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +20.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.096"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "40.45"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.364"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.60%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.9998"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.812.59"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.83%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.226"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001126"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06575"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9990"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.939"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.906.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.227"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.32%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.63"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.11%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.47"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.67%  "
$ws.Range("B28").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C28").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.014.61"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.95%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.405"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.77"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1074"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.12%  "
$ws.Range("E32").Value = "  +0.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.565"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.607"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.82%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06796"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "8.906"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02301"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2139"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.37"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.55%  "
$ws.Range("E40").Value = "  -2.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6196"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9982"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.06"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5861"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.282"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.670"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.949"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.177"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06743"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.17%  "
